# Update the Markov transition-probability matrix on Sheet1 (UIC_A)
# with recomputed values after adding more simulated games.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.198019801980198
$ws.Range("C2").Value = 0.5544554455445545
$ws.Range("J2").Value = 0.01485148514851485
$ws.Range("P2").Value = 0.1386138613861386
$ws.Range("S2").Value = 0.09405940594059406

# Row 3
$ws.Range("B3").Value = 0.01724137931034483
$ws.Range("C3").Value = 0.03448275862068965
$ws.Range("J3").Value = 0.0603448275862069
$ws.Range("P3").Value = 0.646551724137931
$ws.Range("S3").Value = 0.2413793103448276

# Row 4
$ws.Range("J4").Value = 0.05263157894736842
$ws.Range("P4").Value = 0.631578947368421
$ws.Range("S4").Value = 0.3157894736842105

# Row 6
$ws.Range("B6").Value = 0.03589743589743589
$ws.Range("D6").Value = 0.005128205128205128
$ws.Range("F6").Value = 0.07692307692307693
$ws.Range("J6").Value = 0.2564102564102564
$ws.Range("O6").Value = 0.01538461538461539
$ws.Range("Q6").Value = 0.1435897435897436
$ws.Range("R6").Value = 0.09743589743589744
$ws.Range("S6").Value = 0.3692307692307693

# Row 7
$ws.Range("B7").Value = 0.09178743961352658
$ws.Range("D7").Value = 0.02898550724637681
$ws.Range("F7").Value = 0.05314009661835749
$ws.Range("J7").Value = 0.1256038647342995
$ws.Range("O7").Value = 0.02415458937198068
$ws.Range("Q7").Value = 0.1835748792270532
$ws.Range("R7").Value = 0.0821256038647343
$ws.Range("S7").Value = 0.4106280193236715

# Row 8
$ws.Range("B8").Value = 0.06569343065693431
$ws.Range("D8").Value = 0.0145985401459854
$ws.Range("F8").Value = 0.0364963503649635
$ws.Range("J8").Value = 0.09732360097323602
$ws.Range("O8").Value = 0.0170316301703163
$ws.Range("Q8").Value = 0.2043795620437956
$ws.Range("R8").Value = 0.1289537712895377
$ws.Range("S8").Value = 0.4355231143552312

# Row 9
$ws.Range("B9").Value = 0.04761904761904762
$ws.Range("D9").Value = 0.006802721088435374
$ws.Range("F9").Value = 0.06122448979591837
$ws.Range("J9").Value = 0.08843537414965986
$ws.Range("O9").Value = 0.02040816326530612
$ws.Range("Q9").Value = 0.1904761904761905
$ws.Range("R9").Value = 0.1156462585034014
$ws.Range("S9").Value = 0.4693877551020408

# Row 10
$ws.Range("B10").Value = 0.08556611927398444
$ws.Range("D10").Value = 0.02247191011235955
$ws.Range("E10").Value = 0.000864304235090752
$ws.Range("F10").Value = 0.07519446845289542
$ws.Range("J10").Value = 0.09075194468452895
$ws.Range("O10").Value = 0.01815038893690579
$ws.Range("Q10").Value = 0.1961970613656007
$ws.Range("R10").Value = 0.09334485738980121
$ws.Range("S10").Value = 0.4174589455488332

# Row 11
$ws.Range("G11").Value = 0.1314878892733564
$ws.Range("J11").Value = 0.08650519031141868
$ws.Range("K11").Value = 0.1660899653979239
$ws.Range("L11").Value = 0.5986159169550173
$ws.Range("S11").Value = 0.01730103806228374

# Row 12
$ws.Range("G12").Value = 0.7955801104972375
$ws.Range("J12").Value = 0.1325966850828729
$ws.Range("K12").Value = 0.01104972375690608
$ws.Range("L12").Value = 0.03314917127071823
$ws.Range("S12").Value = 0.02762430939226519

# Row 13
$ws.Range("G13").Value = 0.7021276595744681
$ws.Range("J13").Value = 0.2765957446808511
$ws.Range("S13").Value = 0.02127659574468085

# Row 15
$ws.Range("F15").Value = 0.009708737864077669
$ws.Range("H15").Value = 0.1893203883495146
$ws.Range("I15").Value = 0.04854368932038835
$ws.Range("J15").Value = 0.3592233009708738
$ws.Range("K15").Value = 0.04854368932038835
$ws.Range("M15").Value = 0.02427184466019417
$ws.Range("N15").Value = 0.004854368932038835
$ws.Range("O15").Value = 0.06310679611650485
$ws.Range("S15").Value = 0.2524271844660194

# Row 16
$ws.Range("F16").Value = 0.02439024390243903
$ws.Range("H16").Value = 0.1300813008130081
$ws.Range("I16").Value = 0.04065040650406504
$ws.Range("J16").Value = 0.3821138211382114
$ws.Range("K16").Value = 0.1869918699186992
$ws.Range("M16").Value = 0.02439024390243903
$ws.Range("N16").Value = 0.008130081300813009
$ws.Range("O16").Value = 0.08943089430894309
$ws.Range("S16").Value = 0.1138211382113821

# Row 17
$ws.Range("F17").Value = 0.009925558312655087
$ws.Range("H17").Value = 0.173697270471464
$ws.Range("I17").Value = 0.08188585607940446
$ws.Range("J17").Value = 0.424317617866005
$ws.Range("K17").Value = 0.109181141439206
$ws.Range("M17").Value = 0.02233250620347394
$ws.Range("N17").Value = 0.002481389578163772
$ws.Range("O17").Value = 0.07444168734491315
$ws.Range("S17").Value = 0.1017369727047146

# Row 18
$ws.Range("F18").Value = 0.02816901408450704
$ws.Range("H18").Value = 0.1830985915492958
$ws.Range("I18").Value = 0.06572769953051644
$ws.Range("J18").Value = 0.4272300469483568
$ws.Range("K18").Value = 0.1173708920187793
$ws.Range("M18").Value = 0.02816901408450704
$ws.Range("O18").Value = 0.04694835680751173
$ws.Range("S18").Value = 0.1032863849765258

# Row 19
$ws.Range("F19").Value = 0.01813685078318219
$ws.Range("H19").Value = 0.2077493816982688
$ws.Range("I19").Value = 0.07089859851607584
$ws.Range("J19").Value = 0.393239901071723
$ws.Range("K19").Value = 0.1096455070074196
$ws.Range("M19").Value = 0.02143446001648805
$ws.Range("O19").Value = 0.05935696619950536
$ws.Range("S19").Value = 0.1195383347073372
